$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at the top of the data table
# (row 143), pushing the existing rows 143:165 down to 144:166.
$ws.Rows("143:143").Insert()

# Populate the newly inserted row 143 with this week's values.
$ws.Range("A143").Value = 7
$ws.Range("B143").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C143").Value = "Ñuble"
$ws.Range("D143").Value = 44474
$ws.Range("E143").Value = 16
$ws.Range("F143").Value = 100112009
$ws.Range("G143").Value = "Acelga"
$ws.Range("H143").Value = "Sin especificar"
$ws.Range("I143").Value = "Primera"
$ws.Range("J143").Value = 160
$ws.Range("K143").Value = 350
$ws.Range("L143").Value = 400
$ws.Range("M143").Value = 375
$ws.Range("N143").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O143").Value = "Provincia de Diguillín"
$ws.Range("P143").Value = 375
$ws.Range("Q143").Value = 1
$ws.Range("R143").Value = "Hortaliza"
